# Updates cryptos list data (prices / 1h volume changes) to match the
# latest scrape. Leading "'" forces Excel to keep these as text values
# (matching the original inlineStr cell type) instead of auto-converting
# number-like strings (e.g. "68.353.73", "1.00") into numeric/date values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'68.353.73"
$ws.Range("E2").Value = "'  +0.74%  "
$ws.Range("D3").Value = "'2.542.50"
$ws.Range("E3").Value = "'  +0.55%  "
$ws.Range("E4").Value = "'  +0.02%  "
$ws.Range("D5").Value = "'593.16"
$ws.Range("E5").Value = "'  +0.07%  "
$ws.Range("D6").Value = "'175.82"
$ws.Range("E6").Value = "'  -0.17%  "
$ws.Range("E7").Value = "'  -0.02%  "
$ws.Range("E8").Value = "'  -0.76%  "
$ws.Range("D9").Value = "'2.542.82"
$ws.Range("E9").Value = "'  +0.54%  "
$ws.Range("D10").Value = "'0.138"
$ws.Range("E10").Value = "'  -1.71%  "
$ws.Range("E11").Value = "'  +1.72%  "
$ws.Range("D12").Value = "'0.345"
$ws.Range("E12").Value = "'  +0.51%  "
$ws.Range("D13").Value = "'5.03"
$ws.Range("D14").Value = "'26.61"
$ws.Range("E14").Value = "'  -0.80%  "
$ws.Range("D15").Value = "'2.956.90"
$ws.Range("E15").Value = "'  -1.11%  "
$ws.Range("D16").Value = "'0.0000177"
$ws.Range("E16").Value = "'  -0.31%  "
$ws.Range("D17").Value = "'68.385.25"
$ws.Range("E17").Value = "'  +0.94%  "
$ws.Range("B18").Value = "'WrappedEther"
$ws.Range("C18").Value = "'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "'2.693.13"
$ws.Range("E18").Value = "'  +6.35%  "
$ws.Range("B19").Value = "'Chainlink"
$ws.Range("C19").Value = "'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D19").Value = "'11.91"
$ws.Range("E19").Value = "'  +3.87%  "
$ws.Range("B20").Value = "'Uniswap"
$ws.Range("C20").Value = "'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").Value = "'8.06"
$ws.Range("E20").Value = "'  -0.29%  "
$ws.Range("B21").Value = "'Binance-PegBSC-USD"
$ws.Range("C21").Value = "'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D21").Value = "'1.72"
$ws.Range("E21").Value = "'  +72.00%  "
$ws.Range("D22").Value = "'370.52"
$ws.Range("E22").Value = "'  +3.09%  "
$ws.Range("E23").Value = "'  -0.42%  "
$ws.Range("D24").Value = "'4.58"
$ws.Range("D25").Value = "'71.85"
$ws.Range("E25").Value = "'  +1.56%  "
$ws.Range("D27").Value = "'1.91"
$ws.Range("E27").Value = "'  -3.87%  "
$ws.Range("D28").Value = "'9.95"
$ws.Range("E28").Value = "'  -3.13%  "
$ws.Range("D29").Value = "'2.671.42"
$ws.Range("D30").Value = "'" + '0.0' + [char]0x2083 + '0971'
$ws.Range("E30").Value = "'  -1.54%  "
$ws.Range("D31").Value = "'536.10"
$ws.Range("E31").Value = "'  -2.64%  "
$ws.Range("D32").Value = "'8.32"
$ws.Range("E32").Value = "'  +0.46%  "
$ws.Range("E33").Value = "'  -2.69%  "
$ws.Range("D34").Value = "'1.87"
$ws.Range("E34").Value = "'  +0.62%  "
$ws.Range("E35").Value = "'  -0.95%  "
$ws.Range("D36").Value = "'1.00"
$ws.Range("E36").Value = "'  +0.07%  "
$ws.Range("D37").Value = "'158.74"
$ws.Range("E37").Value = "'  +1.05%  "
$ws.Range("E38").Value = "'  -2.11%  "
$ws.Range("D39").Value = "'19.27"
$ws.Range("E39").Value = "'  +2.64%  "
$ws.Range("D40").Value = "'18.63"
$ws.Range("E40").Value = "'  +0.17%  "
$ws.Range("D41").Value = "'5.17"
$ws.Range("E41").Value = "'  +0.00%  "
$ws.Range("D42").Value = "'1.79"
$ws.Range("E42").Value = "'  -1.20%  "
$ws.Range("E43").Value = "'  -1.37%  "
$ws.Range("D44").Value = "'2.54"
$ws.Range("E44").Value = "'  -0.13%  "
$ws.Range("E45").Value = "'  +0.05%  "
$ws.Range("D46").Value = "'39.47"
$ws.Range("E46").Value = "'  -1.26%  "
$ws.Range("D47").Value = "'148.69"
$ws.Range("E47").Value = "'  +0.64%  "
$ws.Range("E48").Value = "'  +0.65%  "
$ws.Range("D49").Value = "'3.73"
$ws.Range("E49").Value = "'  +0.64%  "
$ws.Range("D50").Value = "'0.554"
$ws.Range("E50").Value = "'  -0.98%  "
$ws.Range("E51").Value = "'  +1.84%  "
